$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1631067961165049
$ws.Range("C2").Value = 0.6077669902912621
$ws.Range("J2").Value = 0.01747572815533981
$ws.Range("P2").Value = 0.1203883495145631
$ws.Range("S2").Value = 0.0912621359223301
$ws.Range("B3").Value = 0.006191950464396285
$ws.Range("C3").Value = 0.03715170278637771
$ws.Range("J3").Value = 0.03405572755417956
$ws.Range("P3").Value = 0.7244582043343654
$ws.Range("S3").Value = 0.1981424148606811
$ws.Range("J4").Value = 0.02352941176470588
$ws.Range("P4").Value = 0.6588235294117647
$ws.Range("S4").Value = 0.3176470588235294
$ws.Range("B6").Value = 0.06382978723404255
$ws.Range("D6").Value = 0.009456264775413711
$ws.Range("E6").Value = 0.002364066193853428
$ws.Range("F6").Value = 0.07092198581560284
$ws.Range("J6").Value = 0.2387706855791962
$ws.Range("O6").Value = 0.02127659574468085
$ws.Range("Q6").Value = 0.1654846335697399
$ws.Range("R6").Value = 0.04491725768321513
$ws.Range("S6").Value = 0.3829787234042553
$ws.Range("B7").Value = 0.1313131313131313
$ws.Range("D7").Value = 0.0202020202020202
$ws.Range("F7").Value = 0.05387205387205387
$ws.Range("J7").Value = 0.1212121212121212
$ws.Range("O7").Value = 0.02356902356902357
$ws.Range("Q7").Value = 0.1683501683501684
$ws.Range("R7").Value = 0.09090909090909091
$ws.Range("S7").Value = 0.3905723905723906
$ws.Range("B8").Value = 0.09058402860548272
$ws.Range("D8").Value = 0.02383790226460071
$ws.Range("E8").Value = 0.001191895113230036
$ws.Range("F8").Value = 0.07866507747318235
$ws.Range("J8").Value = 0.08462455303933254
$ws.Range("O8").Value = 0.01430274135876043
$ws.Range("Q8").Value = 0.1728247914183552
$ws.Range("R8").Value = 0.08820023837902265
$ws.Range("S8").Value = 0.4457687723480334
$ws.Range("B9").Value = 0.07932011331444759
$ws.Range("D9").Value = 0.0169971671388102
$ws.Range("E9").Value = 0.0028328611898017
$ws.Range("F9").Value = 0.07365439093484419
$ws.Range("J9").Value = 0.1388101983002833
$ws.Range("O9").Value = 0.0198300283286119
$ws.Range("Q9").Value = 0.2181303116147309
$ws.Range("R9").Value = 0.09915014164305949
$ws.Range("S9").Value = 0.3512747875354107
$ws.Range("B10").Value = 0.1176470588235294
$ws.Range("D10").Value = 0.02380952380952381
$ws.Range("E10").Value = 0.001400560224089636
$ws.Range("F10").Value = 0.07282913165266107
$ws.Range("J10").Value = 0.09570494864612512
$ws.Range("O10").Value = 0.015406162464986
$ws.Range("Q10").Value = 0.2371615312791783
$ws.Range("R10").Value = 0.07889822595704948
$ws.Range("S10").Value = 0.3571428571428572
$ws.Range("G11").Value = 0.1353535353535354
$ws.Range("J11").Value = 0.1131313131313131
$ws.Range("K11").Value = 0.1878787878787879
$ws.Range("L11").Value = 0.5515151515151515
$ws.Range("S11").Value = 0.01212121212121212
$ws.Range("G12").Value = 0.7256317689530686
$ws.Range("J12").Value = 0.1949458483754513
$ws.Range("K12").Value = 0.01083032490974729
$ws.Range("L12").Value = 0.03610108303249097
$ws.Range("S12").Value = 0.03249097472924187
$ws.Range("F13").Value = 0.01639344262295082
$ws.Range("G13").Value = 0.6557377049180327
$ws.Range("J13").Value = 0.2786885245901639
$ws.Range("S13").Value = 0.04918032786885246
$ws.Range("G14").Value = 0.25
$ws.Range("J14").Value = 0.5
$ws.Range("S14").Value = 0.25
$ws.Range("F15").Value = 0.01354401805869074
$ws.Range("H15").Value = 0.1647855530474041
$ws.Range("I15").Value = 0.09255079006772009
$ws.Range("J15").Value = 0.3702031602708803
$ws.Range("K15").Value = 0.06094808126410835
$ws.Range("M15").Value = 0.009029345372460496
$ws.Range("N15").Value = 0.002257336343115124
$ws.Range("O15").Value = 0.04740406320541761
$ws.Range("S15").Value = 0.2392776523702032
$ws.Range("F16").Value = 0.0171919770773639
$ws.Range("H16").Value = 0.2091690544412607
$ws.Range("I16").Value = 0.05444126074498568
$ws.Range("J16").Value = 0.4011461318051576
$ws.Range("K16").Value = 0.1260744985673352
$ws.Range("M16").Value = 0.01146131805157593
$ws.Range("O16").Value = 0.08022922636103152
$ws.Range("S16").Value = 0.1002865329512894
$ws.Range("F17").Value = 0.02375296912114014
$ws.Range("H17").Value = 0.2042755344418052
$ws.Range("I17").Value = 0.09857482185273159
$ws.Range("J17").Value = 0.4026128266033254
$ws.Range("K17").Value = 0.08194774346793349
$ws.Range("M17").Value = 0.01781472684085511
$ws.Range("N17").Value = 0.002375296912114014
$ws.Range("O17").Value = 0.06294536817102138
$ws.Range("S17").Value = 0.1057007125890736
$ws.Range("F18").Value = 0.01904761904761905
$ws.Range("H18").Value = 0.1936507936507937
$ws.Range("I18").Value = 0.09206349206349207
$ws.Range("J18").Value = 0.3841269841269842
$ws.Range("K18").Value = 0.08571428571428572
$ws.Range("M18").Value = 0.02857142857142857
$ws.Range("N18").Value = 0.003174603174603175
$ws.Range("O18").Value = 0.08888888888888889
$ws.Range("S18").Value = 0.1047619047619048
$ws.Range("F19").Value = 0.02209944751381215
$ws.Range("H19").Value = 0.2131675874769798
$ws.Range("I19").Value = 0.08379373848987108
$ws.Range("J19").Value = 0.3646408839779006
$ws.Range("K19").Value = 0.1077348066298343
$ws.Range("M19").Value = 0.0147329650092081
$ws.Range("N19").Value = 0.0009208103130755065
$ws.Range("O19").Value = 0.08517495395948435
$ws.Range("S19").Value = 0.1077348066298343
